$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") values were recorded as "6-3-2007-08" which is off
# by one day due to the way the NBA stats site displayed dates. Correct
# them to the proper ISO-style date string "2008-06-03".
#
# NumberFormat is forced to Text ("@") before the assignment so that
# Excel does not auto-parse the "2008-06-03" literal as a date serial
# value; ClearFormats() afterwards restores the cell's default/general
# formatting while keeping the stored value as text.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    if ($cell.Value2 -eq "6-3-2007-08") {
        $cell.NumberFormat = "@"
        $cell.Value = "2008-06-03"
        $cell.ClearFormats()
    }
}
